$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the style of the existing
# header cells (e.g. G1 "sum") by copying its formatting, then fill in
# H2/H3 with 0 values for the two data rows.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
